{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\n// --- Change 1: add citations after \"OpenCV\" and \"ImageAI\" in the first body paragraph ---\nconst searchResults = context.document.body.search(\"The python libraries OpenCV and ImageAI were used\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  searchResults.items[0].insertText(\n    \"The python libraries OpenCV (Bradski 2000) and ImageAI (Olafenwa and Olafenwa 2018) were used\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// --- Change 2: insert a new Body Text paragraph about quantile regression ---\n// Locate the paragraph that starts the \"Differences among treatments...\" methods\n// paragraph so the new paragraph can be inserted right after it (and before the\n// \"We visually investigated...\" paragraph).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst anchorText = \"Differences among treatments, experimental time points\";\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(anchorText) === 0) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (anchorParagraph) {\n  const newParagraphText =\n    \"Quantile regression analyses were also used to test for significant differences among treatments for larger individuals. This was done to account for potential differences in growth rates from density-dependent feeding effects in each jar that were independent of the treatments. Specifically, individuals in each jar that were located on the periphery of the cluster of oysters may have more access to food as compared to individuals in the center of each cluster. As a result, differences in growth between treatments may be more easily identified by evaluating only individuals with positive growth and only individuals at a higher percentile of growth. Quantile regression models were used to test for significant differences in size (area and weight) for individuals at the 80th percentile for each treatment and on the subset of individuals that showed only positive growth throughout the six-week period. The qt function from the quantreg package was used for all analyses (Koenker 2020). As for the linear models above, separate models were developed to test for differences among treatments within each week and for both week and treatment as predictors.\";\n\n  const newParagraph = anchorParagraph.insertParagraph(newParagraphText, Word.InsertLocation.after);\n  newParagraph.style = \"Body Text\";\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $app / $doc resolve against the live session; the document under\n# edit is $word.ActiveDocument.\n\n$d = $word.ActiveDocument\n\n# --- Change 1: add citations after \"OpenCV\" and \"ImageAI\" in the first body paragraph ---\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n# wdReplaceOne = 1 (replace only the first/next match found)\n$find.Execute(\n    \"The python libraries OpenCV and ImageAI were used\",\n    $false,\n    $false,\n    $false,\n    $false,\n    $false,\n    $true,\n    1,\n    $false,\n    \"The python libraries OpenCV (Bradski 2000) and ImageAI (Olafenwa and Olafenwa 2018) were used\",\n    1\n)\n\n# --- Change 2: insert a new Body Text paragraph about quantile regression ---\n# Find the \"Differences among treatments...\" paragraph so the new paragraph can be\n# inserted right after it (and before the \"We visually investigated...\" paragraph).\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Differences among treatments, experimental time points\")) {\n        $targetPara = $p\n        break\n    }\n}\n\nif ($targetPara -ne $null) {\n    $targetPara.Range.InsertParagraphAfter()\n    $newPara = $targetPara.Next()\n    $newPara.Range.Text = \"Quantile regression analyses were also used to test for significant differences among treatments for larger individuals. This was done to account for potential differences in growth rates from density-dependent feeding effects in each jar that were independent of the treatments. Specifically, individuals in each jar that were located on the periphery of the cluster of oysters may have more access to food as compared to individuals in the center of each cluster. As a result, differences in growth between treatments may be more easily identified by evaluating only individuals with positive growth and only individuals at a higher percentile of growth. Quantile regression models were used to test for significant differences in size (area and weight) for individuals at the 80th percentile for each treatment and on the subset of individuals that showed only positive growth throughout the six-week period. The qt function from the quantreg package was used for all analyses (Koenker 2020). As for the linear models above, separate models were developed to test for differences among treatments within each week and for both week and treatment as predictors.\"\n    $newPara.Range.Style = \"Body Text\"\n}\n"}
